# Insert a new weekly price record for "Poroto verde" (Femacal de La Calera)
# immediately before the current row 161, shifting the existing rows 161-215
# down to 162-216. The sheet's dimension grows from A1:R215 to A1:R216.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 161; everything from the old row 161
# downward (through row 215) moves down one row, i.e. old row 161 becomes
# new row 162, ..., old row 215 becomes new row 216.
$ws.Rows.Item(161).Insert()

# Populate the newly inserted row 161 with the new record's data.
$ws.Range("A161").Value = 3
$ws.Range("B161").Value = "Femacal de La Calera"
$ws.Range("C161").Value = "Coquimbo"
$ws.Range("D161").Value2 = 44468
$ws.Range("E161").Value = 5
$ws.Range("F161").Value = 100112031
$ws.Range("G161").Value = "Poroto verde"
$ws.Range("H161").Value = "Magnum"
$ws.Range("I161").Value = "Primera"
$ws.Range("J161").Value = 65
$ws.Range("K161").Value = 34000
$ws.Range("L161").Value = 35000
$ws.Range("M161").Value = 34462
$ws.Range("N161").Value = "`$/malla 25 kilos"
$ws.Range("O161").Value = "Región de Arica y Parinacota"
$ws.Range("P161").Value = 1378
$ws.Range("Q161").Value = 25
$ws.Range("R161").Value = "Hortaliza"
